$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('Associate - Data Engineer - Global Master Data', 'Eli Lilly', 'Indianapolis, IN, US USA', 15.6, 'Data Scientist, RAG, Glue, Redshift, Data Lake, CI/CD, Git, Databricks, Redshift, Kafka', '2026-02-23', 'https://www.indeed.com/viewjob?jk=326b602342da8867')
    ,('IT Data Platform Engineer', 'National Retail Systems, Inc', 'Lyndhurst, NJ, US USA', 15.6, 'Data Scientist, RAG, S3, Glue, Athena, CI/CD, Git, Snowflake, Power BI, Python', '2026-02-23', 'https://www.indeed.com/viewjob?jk=4da5376ea9ed63ff')
    ,('Data Scientist III', 'Fanatics', 'US USA', 15.6, 'Data Scientist, Redshift, BigQuery, MLflow, Git, Snowflake, Databricks, BigQuery, Redshift, Python', '2026-02-23', 'https://www.indeed.com/viewjob?jk=a21ec81d7d224fd0')
    ,('Full Stack Magento Developer – Remote (Canada)', 'Jarvis Recruitment', 'Remote, US USA', 14.4, 'RAG, Prompt Engineering, Docker, CI/CD, Jenkins, GitHub Actions, Git, MySQL, SQL, R', '2026-02-23', 'https://www.indeed.com/viewjob?jk=5e6328c1040c5967')
    ,('Analytics Engineer', 'RollKall', 'Irving, TX, US USA', 11.1, 'RAG, Dataflow, Git, Tableau, Power BI, Python, SQL, R, Scala, Optimization', '2026-02-23', 'https://www.indeed.com/viewjob?jk=68db2c37ad6b9984')
    ,('Machine Learning Engineer', 'Transflo', 'Remote, US USA', 11.1, 'Data Scientist, Machine Learning Engineer, Generative AI, RAG, Prompt Engineering, MLflow, CI/CD, R, Scala, Optimization', '2026-02-23', 'https://www.indeed.com/viewjob?jk=366c91c1a2a15b5c')
    ,('Machine Learning Engineer 50', 'Adobe', 'San Jose, CA, US USA', 11.1, 'Machine Learning Engineer, Generative AI, RAG, TensorFlow, PyTorch, Docker, Power BI, Python, R, Scala', '2026-02-23', 'https://www.indeed.com/viewjob?jk=652c1e250e125af4')
    ,('Sr. Cloud Engineer', 'Transflo', 'Remote, US USA', 10, 'S3, CI/CD, Jenkins, GitHub Actions, Terraform, Git, Python, R, Optimization', '2026-02-23', 'https://www.indeed.com/viewjob?jk=ce6390023da67737')
    ,('Machine Learning Engineer', 'Adobe', 'San Jose, CA, US USA', 10, 'Machine Learning Engineer, Generative AI, RAG, TensorFlow, PyTorch, CI/CD, Python, R, Scala', '2026-02-23', 'https://www.indeed.com/viewjob?jk=1c4f153b447e2441')
    ,('Software Engineer II', 'Blue Origin', 'Seattle, WA, US USA', 10, 'RAG, CI/CD, Git, NoSQL, Python, SQL, R, Java, Scala', '2026-02-23', 'https://www.indeed.com/viewjob?jk=3b3a27f0f05bdbda')
    ,('Senior Salesforce Tableau Developer - Sales Focused', 'Ferguson', 'Remote, US USA', 10, 'RAG, Synapse, Git, Databricks, Tableau, Power BI, SQL, R, Scala', '2026-02-23', 'https://www.indeed.com/viewjob?jk=e85774318b17a840')
    ,('Specialist - Architecture', 'LTIMindtree', 'New York, NY, US USA', 10, 'Generative AI, RAG, LLaMA, Pinecone, Prompt Engineering, Python, R, Java, Hypothesis Testing', '2026-02-23', 'https://www.indeed.com/viewjob?jk=8dfb9f688fa2bf88')
    ,('Specialist - Architecture', 'LTIMindtree', 'New York, NY, US USA', 10, 'Generative AI, RAG, LLaMA, Pinecone, Prompt Engineering, Python, R, Java, Hypothesis Testing', '2026-02-23', 'https://www.indeed.com/viewjob?jk=b882e08a5c9eb997')
    ,('Specialist - Architecture', 'LTIMindtree', 'New York, NY, US USA', 10, 'Generative AI, RAG, LLaMA, Pinecone, Prompt Engineering, Python, R, Java, Hypothesis Testing', '2026-02-23', 'https://www.indeed.com/viewjob?jk=e18794724c7fb0b4')
    ,('Specialist - Architecture', 'LTIMindtree', 'New York, NY, US USA', 10, 'Generative AI, RAG, LLaMA, Pinecone, Prompt Engineering, Python, R, Java, Hypothesis Testing', '2026-02-23', 'https://www.indeed.com/viewjob?jk=baa2a66d4c015d1f')
    ,('Specialist - Architecture', 'LTIMindtree', 'New York, NY, US USA', 10, 'Generative AI, RAG, LLaMA, Pinecone, Prompt Engineering, Python, R, Java, Hypothesis Testing', '2026-02-23', 'https://www.indeed.com/viewjob?jk=ba86dc84f904b800')
    ,('Bioinformatics Software Engineer', 'Harvard University', 'Boston, MA, US USA', 10, 'RAG, S3, EC2, Docker, CI/CD, Git, Python, R, Scala', '2026-02-23', 'https://www.indeed.com/viewjob?jk=8d28ab892406e1df')
)

# Keep the "Posted At" column as plain text (not auto-converted to a date serial)
$ws.Range("F2:F18").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
